$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3 for the new "instructor" user, shifting all
# existing W*/D* rows down by one.
$ws.Rows(3).Insert()

# New instructor account (row 3)
$ws.Range("A3").Value = "D010"
$ws.Range("B3").Value = "passInstructor"

# Update the admin account's password (row 2) to include special characters
$ws.Range("B2").Value = "passadmin_<>?"

# New "Role" column (C) - header + roles for the two accounts that have one
$ws.Range("C1").Value = "Role"
$ws.Range("C1").Font.Bold = $true
$ws.Range("C2").Value = "admin"
$ws.Range("C3").Value = "instructor"

# Give column C a sensible width, matching the other data columns
$ws.Columns("C").ColumnWidth = 14.14

# Reset scroll position / selection to match the edited area
$ws.Range("C9").Select()
